# Apply cryptos list price/volume update (commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.933.52"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "'1.639.35"
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.62%  "
$ws.Range("D5").Value = "'214.73"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").Value = "'0.5078"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("D8").Value = "'0.2576"
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").Value = "'0.06355"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").Value = "'19.80"
$ws.Range("E10").Value = "  +1.17%  "
$ws.Range("D11").Value = "'0.07727"
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("D12").Value = "'4.293"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").Value = "'1.642.85"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("D14").Value = "'0.5463"
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").Value = "'0.0₅7733"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("D16").Value = "'64.23"
$ws.Range("E16").Value = "  -0.82%  "
$ws.Range("D17").Value = "'25.966.30"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("D20").Value = "'195.88"
$ws.Range("D21").Value = "'9.945"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").Value = "'6.136"
$ws.Range("E22").Value = "  +1.54%  "
$ws.Range("D23").Value = "'1.003"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").Value = "'1.892"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").Value = "'142.69"
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("E26").Value = "  +10.92%  "
$ws.Range("D27").Value = "'6.856"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("D28").Value = "'15.58"
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("D29").Value = "'1.239"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "'0.04888"
$ws.Range("E30").Value = "  -2.88%  "
$ws.Range("D31").Value = "'3.263"
$ws.Range("D32").Value = "'3.203"
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").Value = "'2.376"
$ws.Range("E34").Value = "  +0.43%  "
$ws.Range("D35").Value = "'0.9159"
$ws.Range("E35").Value = "  +2.29%  "
$ws.Range("D36").Value = "'2.569"
$ws.Range("E36").Value = "  -0.87%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.5532"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "'1.131.29"
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("D39").Value = "'0.01567"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("D40").Value = "'1.001"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").Value = "'5.589"
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("D42").Value = "'0.8036"
$ws.Range("E42").Value = "  -1.69%  "
$ws.Range("D43").Value = "'98.61"
$ws.Range("E43").Value = "  -1.42%  "
$ws.Range("D44").Value = "'0.0₈120"
$ws.Range("E44").Value = "  -9.16%  "
$ws.Range("D45").Value = "'1.777.08"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("D46").Value = "'0.4521"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").Value = "'55.22"
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("D48").Value = "'0.9984"
$ws.Range("E48").Value = "  -0.84%  "
$ws.Range("D49").Value = "'0.05185"
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("D50").Value = "'7.482"
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("E51").Value = "  -0.43%  "
